$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# Text-like columns: use a leading apostrophe (quote-prefix) so Excel keeps
# these as literal text instead of auto-coercing to a date / number.
$ws.Cells.Item($row, 1).Value = "'2024-01-05"
$ws.Cells.Item($row, 2).Value = "15:09:18"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "'00"

# Numeric columns
$ws.Cells.Item($row, 5).Value = 140627
$ws.Cells.Item($row, 6).Value = 142923
$ws.Cells.Item($row, 7).Value = 172160
$ws.Cells.Item($row, 8).Value = 147050
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 117968
$ws.Cells.Item($row, 11).Value = 224447
$ws.Cells.Item($row, 12).Value = 248611
$ws.Cells.Item($row, 13).Value = 184701
$ws.Cells.Item($row, 14).Value = 110132
$ws.Cells.Item($row, 15).Value = 40431
$ws.Cells.Item($row, 16).Value = 30810
$ws.Cells.Item($row, 17).Value = 72384
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41492
$ws.Cells.Item($row, 20).Value = -1
